# Insert a new row at position 210 (this shifts existing rows 210:303 down to 211:304,
# matching the diff where dimension grows from A1:T303 to A1:T304 and every row from
# 211..304 now holds what used to be in 210..303).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("210:210").Insert()

# Populate the newly inserted row 210 with its data.
$ws.Cells.Item(210, 1).Value = 11
$ws.Cells.Item(210, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(210, 3).Value = "Bíobío"
$ws.Cells.Item(210, 4).Value = 44755
$ws.Cells.Item(210, 5).Value = 8
$ws.Cells.Item(210, 6).Value = "Fruta"
$ws.Cells.Item(210, 7).Value = 100102
$ws.Cells.Item(210, 8).Value = "Cítricos"
$ws.Cells.Item(210, 9).Value = 100102005
$ws.Cells.Item(210, 10).Value = "Naranja"
$ws.Cells.Item(210, 11).Value = "Fukumoto"
$ws.Cells.Item(210, 12).Value = "Primera"
$ws.Cells.Item(210, 13).Value = 200
$ws.Cells.Item(210, 14).Value = 7000
$ws.Cells.Item(210, 15).Value = 8000
$ws.Cells.Item(210, 16).Value = 7500
$ws.Cells.Item(210, 17).Value = "`$/caja 15 kilos granel"
$ws.Cells.Item(210, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(210, 19).Value = 500
$ws.Cells.Item(210, 20).Value = 15
